$wb = $excel.ActiveWorkbook

# The workbook has two sheets with identical content: "展览" and "全部类型".
# Both need the same updates applied (matching the diff, which touches both
# sheet1 and sheet4 identically).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: 最低票价 (lowest ticket price) 68 -> 78
    $ws.Range("G2").Value = 78

    # Row 3: 想去人数 (want-to-go count) 1263 -> 1268
    $ws.Range("F3").Value = 1268

    # Row 4: 想去人数 1542 -> 1556
    $ws.Range("F4").Value = 1556

    # Row 6: 想去人数 6165 -> 6172
    $ws.Range("F6").Value = 6172
}
